# Applies the diff to UserStorytable_ModifiedAlertsPart_LiyuanQin.docx
#
# wdReplace constants used below:
#   0 = wdReplaceNone   (no replace, just find)
#   1 = wdReplaceOne    (replace first match only)
#   2 = wdReplaceAll    (replace all matches in range)
# wdFindWrap:
#   1 = wdFindContinue

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# ---------------------------------------------------------------------
# Row 19 (Alerts summary story), column 2, 2nd paragraph:
# "As a User, I want to get specific alerts for next arriving buses and
#  general alerts for delayed/cancelled buses and, so that I can<nbsp>spend
#  as little time as possible waiting for the bus and<nbsp>change my plan
#  according to the real-time route status"
# ---------------------------------------------------------------------
$cell19 = $t.Cell(19, 2)

# 1) Remove the _GoBack bookmark that sits right after "...want to get ".
#    Find/Replace across a span that includes the bookmark position drops it,
#    same as real Word collapsing an empty bookmark touched by an edit.
$r = $cell19.Range
$old = "As a User, I want to get specific alerts for "
$null = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $old, 2)

# 2) "delayed/cancelled buses and" + ", so " -> "delayed/cancelled buses " + "so "
$r = $cell19.Range
$old = "delayed/cancelled buses and, so "
$new = "delayed/cancelled buses so "
$null = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

# 3) "can<nbsp>spend as little time as possible" -> "can spend less time"
$r = $cell19.Range
$old = [char]0x0063 + "an" + [char]0x00A0 + "spend as little time as possible"
$new = "can spend less time"
$null = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

Write-Output "Row 19 done: $($cell19.Range.Text)"
